$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A, shifting the existing size/type/experience headers
# (and the column they occupy) one column to the right.
$ws.Columns.Item(1).Insert()

# New first column: hall_id codes (entered with a leading quote so Excel treats
# them as quote-prefixed text, e.g. 'H01' -> stored text 'H01' with quotePrefix style)
$ws.Range("A2").Value = "''H01'"
$ws.Range("A3").Value = "''H02'"
$ws.Range("A4").Value = "''H05'"
$ws.Range("A5").Value = "''H07'"
$ws.Range("A6").Value = "''H09'"
$ws.Range("A7").Value = "''H12'"
$ws.Range("A8").Value = "''H13'"
$ws.Range("A9").Value = "''H16'"
$ws.Range("A10").Value = "''H17'"
$ws.Range("A11").Value = "''H18'"

# New header label for the inserted column
$ws.Range("A1").Value = "hall_id"

# Remaining columns: size (B), type (C), experience (D)
$ws.Range("B2").Value = "''4500 sq feet'"
$ws.Range("C2").Value = "''Dome'"
$ws.Range("D2").Value = "''2D'"

$ws.Range("B3").Value = "''5000 sq feet'"
$ws.Range("C3").Value = "''Flat'"
$ws.Range("D3").Value = "''2D'"

$ws.Range("B4").Value = "''3500 sq feet'"
$ws.Range("C4").Value = "''Curve'"
$ws.Range("D4").Value = "''3D'"

$ws.Range("B5").Value = "''2500 sq feet'"
$ws.Range("C5").Value = "''Flat'"
$ws.Range("D5").Value = "''3D'"

$ws.Range("B6").Value = "''4500 sq feet'"
$ws.Range("C6").Value = "''Dome'"
$ws.Range("D6").Value = "''4D'"

$ws.Range("B7").Value = "''5000 sq feet'"
$ws.Range("C7").Value = "''Dome'"
$ws.Range("D7").Value = "''3D'"

$ws.Range("B8").Value = "''4500 sq feet'"
$ws.Range("C8").Value = "''Flat'"
$ws.Range("D8").Value = "''2D'"

$ws.Range("B9").Value = "''3500 sq feet'"
$ws.Range("C9").Value = "''Curve'"
$ws.Range("D9").Value = "''2D'"

$ws.Range("B10").Value = "''4500 sq feet'"
$ws.Range("C10").Value = "''Flat'"
$ws.Range("D10").Value = "''3D'"

$ws.Range("B11").Value = "''3500 sq feet'"
$ws.Range("C11").Value = "''Curve'"
$ws.Range("D11").Value = "''2D'"

# Column widths (auto-fit to the new content, matching the saved widths)
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).AutoFit()

# Sheet view / selection
$ws.Range("J6").Select()
